# "adicionei o caso de uso 'eliminar promoção'" (não completo)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1) Fix the wording of the existing "Lançar promoção" objective (row 14, col C)
$ws.Range("C14").Value = "O objetivo é torna-la ativa"

# 2) Add the new use case row (row 20): "Eliminar Promoção"
$ws.Range("B20").Value = "Eliminar Promoção"
$ws.Range("C20").Value = "O objetivo é torna-la inativa"

# 3) Extend the "Gestor de Marketing" actor merge from A13:A19 to A13:A20
$ws.Range("A13:A19").UnMerge()
$ws.Range("A13:A20").Merge()

# Re-apply the plain border style (matching the rest of column A) since Merge()
# stamps its own outline-border formatting on the merged range.
$ws.Range("A6").Copy()
$ws.Range("A13:A20").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 4) Give the new row's B/C cells the same bordered style as the rest of the table
$ws.Range("B19:C19").Copy()
$ws.Range("B20:C20").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 5) Restore selection similar to the authored file
$ws.Range("G20").Select()
